$d = $word.ActiveDocument

# --- Edit 1 --------------------------------------------------------------
# "...жилого помещения у лиц..." -> "...жилого помещения (${JP_TYPE}) у лиц..."
# with ${JP_TYPE} placed in its own bold run.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    'жилого помещения у лиц',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'жилого помещения (${JP_TYPE}) у лиц', 2)

$rng1b = $d.Content
$found1b = $rng1b.Find.Execute(
    '${JP_TYPE}',
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1b) {
    $rng1b.Bold = $true
}

# --- Edit 2 --------------------------------------------------------------
# "...я и члены моей семьи не участвовали." ->
# "...я и члены моей семьи ${IS_PARTICIPATE}." with ${IS_PARTICIPATE} bold.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    'не участвовали.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    '${IS_PARTICIPATE}.', 2)

$rng2b = $d.Content
$found2b = $rng2b.Find.Execute(
    '${IS_PARTICIPATE}',
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2b) {
    $rng2b.Bold = $true
}
